$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price strings that look like plain numbers (e.g. "245.65").
# Assigning them directly would make Excel auto-convert to a numeric cell, losing the
# exact text representation (trailing zeros, multi-dot "thousands" notation, etc).
# Prefixing with a single quote forces a text entry; re-applying the Normal style
# afterwards clears the transient quote-prefix formatting so the cell style is left
# exactly as it was (no explicit style index), matching plain inline/shared text cells.

$ws.Range("D2").Value = "42.353.30"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.237.73"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'245.65"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'0.621"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("D7").Value = "'74.25"
$ws.Range("E7").Value = "  -3.82%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").Value = "'43.58"
$ws.Range("E10").Value = "  +4.89%  "
$ws.Range("D11").Value = "'0.0962"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "'7.12"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "'14.46"
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "2.258.46"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "42.250.87"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  +11.56%  "
$ws.Range("D19").Value = "'6.17"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").Value = "'72.05"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "'10.34"
$ws.Range("E21").Value = "  +40.35%  "
$ws.Range("D22").Value = "'231.58"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  -4.68%  "
$ws.Range("D24").Value = "'11.74"
$ws.Range("E24").Value = "  +4.21%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").Value = "'2.31"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "'166.78"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("D30").Value = "'20.94"
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").Value = "'5.91"
$ws.Range("E31").Value = "  +20.06%  "
$ws.Range("D32").Value = "'0.0811"
$ws.Range("E32").Value = "  -1.98%  "
$ws.Range("D33").Value = "'0.119"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").Value = "'29.98"
$ws.Range("E34").Value = "  -9.86%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").Value = "'4.52"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("D38").Value = "'13.25"
$ws.Range("E38").Value = "  -6.89%  "
$ws.Range("D39").Value = "'2.17"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "'5.66"
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("D41").Value = "'63.44"
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("D44").Value = "'105.69"
$ws.Range("E44").Value = "  -7.05%  "
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("E48").Value = "  +3.56%  "
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "'2.73"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("E51").Value = "  -1.98%  "

# Reset style on the forced-text cells so no stray quotePrefix style sticks around.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D50").Style = "Normal"
